# Add log progress rows for the remaining strings that still need translation review.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Hello", "Need review"),
    @("House", "Need review"),
    @("Password", "Need review"),
    @("Run if you like", "Need review"),
    @("Wood", "Need review")
)

$r = 4
foreach ($entry in $rows) {
    $name = $entry[0]
    $status = $entry[1]
    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $name
    $ws.Cells.Item($r, 4).Value = $status
    $r++
}
